# Apply cryptos list update (price/volume refresh + two row re-orderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.751.99"
Set-TextCell "E2" "  +1.45%  "

Set-TextCell "D3" "3.819.31"
Set-TextCell "E3" "  +0.28%  "

Set-TextCell "D4" "0.998"
Set-TextCell "E4" "  -0.01%  "

Set-TextCell "D5" "613.28"

Set-TextCell "D6" "164.88"
Set-TextCell "E6" "  -0.70%  "

Set-TextCell "D7" "3.817.31"
Set-TextCell "E7" "  +0.34%  "

Set-TextCell "E8" "  +0.02%  "

Set-TextCell "E9" "  -0.05%  "

Set-TextCell "E10" "  +0.89%  "

Set-TextCell "D11" "0.451"
Set-TextCell "E11" "  -0.33%  "

Set-TextCell "D12" "6.69"
Set-TextCell "E12" "  +5.40%  "

Set-TextCell "D13" "0.0000249"
Set-TextCell "E13" "  -0.47%  "

Set-TextCell "D14" "35.49"
Set-TextCell "E14" "  -1.30%  "

Set-TextCell "D15" "4.459.95"
Set-TextCell "E15" "  +0.23%  "

Set-TextCell "D16" "3.780.38"
Set-TextCell "E16" "  -1.02%  "

Set-TextCell "D17" "68.689.00"
Set-TextCell "E17" "  +1.31%  "

Set-TextCell "E18" "  -1.29%  "

Set-TextCell "D19" "7.11"
Set-TextCell "E19" "  +0.48%  "

Set-TextCell "E20" "  -0.20%  "

Set-TextCell "D21" "463.96"
Set-TextCell "E21" "  -0.03%  "

Set-TextCell "D22" "9.66"
Set-TextCell "E22" "  -1.70%  "

Set-TextCell "E23" "  -0.03%  "

Set-TextCell "D24" "0.0000151"
Set-TextCell "E24" "  +3.60%  "

Set-TextCell "D25" "83.85"
Set-TextCell "E25" "  +0.63%  "

Set-TextCell "D26" "12.04"
Set-TextCell "E26" "  -0.63%  "

Set-TextCell "E27" "  +0.07%  "

Set-TextCell "E28" "  +0.04%  "

Set-TextCell "D29" "9.99"
Set-TextCell "E29" "  -0.20%  "

Set-TextCell "D30" "3.965.29"
Set-TextCell "E30" "  +0.17%  "

Set-TextCell "D31" "2.63"
Set-TextCell "E31" "  -5.40%  "

Set-TextCell "E32" "  +0.20%  "

Set-TextCell "E33" "  -2.18%  "

Set-TextCell "D34" "29.03"
Set-TextCell "E34" "  -1.30%  "

Set-TextCell "E35" "  -0.03%  "

Set-TextCell "E36" "  +0.05%  "

Set-TextCell "E37" "  +1.74%  "

Set-TextCell "E38" "  +6.14%  "

Set-TextCell "E39" "  +1.62%  "

Set-TextCell "D40" "0.982"
Set-TextCell "E40" "  -1.34%  "

Set-TextCell "D41" "3.17"

Set-TextCell "D42" "0.999"
Set-TextCell "E42" "  -0.05%  "

Set-TextCell "B44" "Monero"
Set-TextCell "C44" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D44" "154.10"
Set-TextCell "E44" "  +1.62%  "

Set-TextCell "B45" "TheGraph"
Set-TextCell "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D45" "0.298"
Set-TextCell "E45" "  -0.29%  "

Set-TextCell "D46" "46.59"
Set-TextCell "E46" "  -2.43%  "

Set-TextCell "B47" "Arweave"
Set-TextCell "C47" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextCell "D47" "42.66"
Set-TextCell "E47" "  -4.36%  "

Set-TextCell "B48" "ONDO"
Set-TextCell "C48" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell "D48" "1.39"
Set-TextCell "E48" "  +0.14%  "

Set-TextCell "D49" "8.37"
Set-TextCell "E49" "  +0.33%  "

Set-TextCell "E50" "  +1.76%  "

Set-TextCell "D51" "378.38"
Set-TextCell "E51" "  -2.94%  "
